$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 59 ("Profile58" test case), modeled on row 58 ---

# Copy row 58's cell formatting down to row 59 (A:E), then fix D59 up to
# reuse the plain bordered style (same as E58) instead of the accent style
# that D58 happens to carry.
$ws.Range("A58:E58").Copy()
$ws.Range("A59:E59").PasteSpecial(-4122)
$ws.Range("E58").Copy()
$ws.Range("D59").PasteSpecial(-4122)

# Row 59, like row 58 (and the other multi-line rows), is a tall row.
$ws.Range("A59:E59").RowHeight = 30

# Cell values (new shared strings get created as needed).
$ws.Range("A59").Value = "Profile58"
$ws.Range("B59").Value = "OPQA-2105|OPQA-2103"
$ws.Range("C59").Value = "Verify that user has the ability to update his Name from the profile modal.|Verify that profile modal displays the following fields pre-populated with values from the user's profile:Name (required field),Title/Role,Country,Skills and Interests (Topics)"
$ws.Range("D59").Value = "Y"

# --- Scroll / selection bookkeeping that moved along with the new row ---
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C55").Select()
